# Generate Report for Handback
#
# Marks the two source files as "handed back" (in sync with en-US) on the
# Overview sheet, and on each per-language sheet (zh-cn, de-de) records the
# latest handback: the target (source) file that was sent out, the handback
# (translated) file that came back, and the handback timestamp.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status columns for both rows
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

# Row 2 — 02f102eb-82c7-4ed2-a3e5-de76edd55c7d
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/73c5df500727be70b3c8b47d82e6f4ac9ef7b7b1/e2e/02f102eb-82c7-4ed2-a3e5-de76edd55c7d.md",
    "",
    "",
    "02f102eb-82c7-4ed2-a3e5-de76edd55c7d.md"
) | Out-Null
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75398ef7f80631156069cd4e9b69b277281f2fb3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/02f102eb-82c7-4ed2-a3e5-de76edd55c7d.34832b11ad3cc89d85c9623070d8bb4761b0ba1e.zh-cn.xlf",
    "",
    "",
    "02f102eb-82c7-4ed2-a3e5-de76edd55c7d.34832b11ad3cc89d85c9623070d8bb4761b0ba1e.zh-cn.xlf"
) | Out-Null
$zhcn.Range("H2").Value = "2016-03-25 12:34:17"

# Row 3 — 57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/73c5df500727be70b3c8b47d82e6f4ac9ef7b7b1/e2e/57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.md",
    "",
    "",
    "57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.md"
) | Out-Null
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/75398ef7f80631156069cd4e9b69b277281f2fb3/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.dcab58c23dd06f9dbba17eee240a705787c4ec11.zh-cn.xlf",
    "",
    "",
    "57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.dcab58c23dd06f9dbba17eee240a705787c4ec11.zh-cn.xlf"
) | Out-Null
$zhcn.Range("H3").Value = "2016-03-25 12:34:17"

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $zhcn.Range($addr).Style = "HyperLink"
}

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# Row 2 — 02f102eb-82c7-4ed2-a3e5-de76edd55c7d
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/73c5df500727be70b3c8b47d82e6f4ac9ef7b7b1/e2e/02f102eb-82c7-4ed2-a3e5-de76edd55c7d.md",
    "",
    "",
    "02f102eb-82c7-4ed2-a3e5-de76edd55c7d.md"
) | Out-Null
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cd84f084c1878fd38253590b0361fceb6a1dc29c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/high/02f102eb-82c7-4ed2-a3e5-de76edd55c7d.34832b11ad3cc89d85c9623070d8bb4761b0ba1e.de-de.xlf",
    "",
    "",
    "02f102eb-82c7-4ed2-a3e5-de76edd55c7d.34832b11ad3cc89d85c9623070d8bb4761b0ba1e.de-de.xlf"
) | Out-Null
$dede.Range("H2").Value = "2016-03-25 12:34:27"

# Row 3 — 57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/73c5df500727be70b3c8b47d82e6f4ac9ef7b7b1/e2e/57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.md",
    "",
    "",
    "57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.md"
) | Out-Null
$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/cd84f084c1878fd38253590b0361fceb6a1dc29c/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/high/57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.dcab58c23dd06f9dbba17eee240a705787c4ec11.de-de.xlf",
    "",
    "",
    "57d9f5a6-57a3-4715-bfd3-e5b8509ffc1d.dcab58c23dd06f9dbba17eee240a705787c4ec11.de-de.xlf"
) | Out-Null
$dede.Range("H3").Value = "2016-03-25 12:34:27"

foreach ($addr in @("F2", "G2", "F3", "G3")) {
    $dede.Range($addr).Style = "HyperLink"
}
